$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.116667747497559
$ws.Range("B1").Value = 2.619241952896118
$ws.Range("C1").Value = 2.746314764022827
$ws.Range("D1").Value = 3.084755659103394
$ws.Range("E1").Value = 0.797762393951416
